$d = $word.ActiveDocument

$null = $d.Paragraphs(20).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>disableSprite</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)]</w:t></w:r></w:p>')
$null = $d.Paragraphs(43).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">[Along the way, the human was traumatized by a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vegetoid</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and killed it.]</w:t></w:r></w:p>')
$null = $d.Paragraphs(71).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">Why try to make friends with that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Froggit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?</w:t></w:r></w:p>')
$null = $d.Paragraphs(81).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">[Near </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Toriel’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> house, it was decided—the duo would go back and kill every monster in the Ruins.]</w:t></w:r></w:p>')
$null = $d.Paragraphs(89).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:proofErr w:type="spellStart"/><w:r><w:t>heya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')
$null = $d.Paragraphs(92).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">so, I’ve got a question for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ya</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')
$null = $d.Paragraphs(95).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">heh </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>heh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>heh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>heh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>…</w:t></w:r></w:p>')
$null = $d.Paragraphs(98).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">do you </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wanna</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> have a bad time?</w:t></w:r></w:p>')
$null = $d.Paragraphs(99).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:proofErr w:type="spellStart"/><w:r><w:t>‘cause</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> if you take another step forward…</w:t></w:r></w:p>')
$null = $d.Paragraphs(101).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:proofErr w:type="spellStart"/><w:r><w:t>welp</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$null = $d.Paragraphs(103).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">this is why </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> never make promises.</w:t></w:r></w:p>')
$null = $d.Paragraphs(109).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">whoa, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>buddo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>')
$null = $d.Paragraphs(110).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:proofErr w:type="spellStart"/><w:r><w:t>uhhhhh</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$null = $d.Paragraphs(116).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">[With no more LOVE to corrupt them, Frisk went on their mission to save everyone—every </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>whimsun</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and every </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vulkin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">; every </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>madgick</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and every </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>froggit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.]</w:t></w:r></w:p>')
